$d = $word.ActiveDocument

$replacements = @(
    @{old="51×34=1734"; new="84×76=6384"},
    @{old="42×59=2478"; new="11×83=913"},
    @{old="51×18=918"; new="42×23=966"},
    @{old="69×36=2484"; new="73×81=5913"},
    @{old="23×76=1748"; new="45×84=3780"},
    @{old="86×67=5762"; new="73×87=6351"},
    @{old="89×96=8544"; new="18×14=252"},
    @{old="90×58=5220"; new="46×71=3266"},
    @{old="58×17=986"; new="95×42=3990"},
    @{old="72×30=2160"; new="69×96=6624"},
    @{old="34×14=476"; new="11×62=682"},
    @{old="15×92=1380"; new="94×39=3666"},
    @{old="69×76=5244"; new="23×60=1380"},
    @{old="62×65=4030"; new="77×17=1309"},
    @{old="73×44=3212"; new="19×65=1235"},
    @{old="18×20=360"; new="91×97=8827"},
    @{old="46×80=3680"; new="90×75=6750"},
    @{old="92×46=4232"; new="56×82=4592"},
    @{old="60×24=1440"; new="75×97=7275"},
    @{old="27×66=1782"; new="18×83=1494"},
    @{old="57×70=3990"; new="59×14=826"},
    @{old="26×44=1144"; new="52×54=2808"},
    @{old="53×67=3551"; new="16×96=1536"},
    @{old="32×75=2400"; new="85×53=4505"},
    @{old="99×97=9603"; new="49×22=1078"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
